$wb = $excel.ActiveWorkbook
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "cudaPython"

$ws1 = $wb.Worksheets.Item("GreenIteration")
$r = $ws1.Range("A3:B13")
$r.HorizontalAlignment = -4108
$r.WrapText = $true
